$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3560.6365
$ws.Range("I137").Value = 917.2593000000001
$ws.Range("K137").Value = 2751.7779
$ws.Range("M137").Value = -201.7779
$ws.Range("H141").Value = 3662.1155
$ws.Range("I141").Value = 1375.9375
$ws.Range("J141").Value = 7320
$ws.Range("K141").Value = 4127.8125
$ws.Range("L141").Value = 21960
$ws.Range("M141").Value = 1052.1875
$ws.Range("N141").Value = -32320

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6645.8867
$ws.Range("I32").Value = 6244.7593
$ws.Range("K32").Value = 6244.7593
$ws.Range("M32").Value = -5957.7593
$ws.Range("H61").Value = 1207.2703
$ws.Range("I61").Value = 919.8148
$ws.Range("J61").Value = 1983.4
$ws.Range("K61").Value = 919.8148
$ws.Range("L61").Value = 1983.4
$ws.Range("M61").Value = -707.8148
$ws.Range("N61").Value = -2407.4
$ws.Range("H74").Value = 1657.6136
$ws.Range("I74").Value = 1466.9429
$ws.Range("J74").Value = 2399.111
$ws.Range("K74").Value = 1466.9429
$ws.Range("L74").Value = 2399.111
$ws.Range("M74").Value = -592.9429
$ws.Range("N74").Value = -4147.111
$ws.Range("H77").Value = 1657.6136
$ws.Range("I77").Value = 1466.9429
$ws.Range("J77").Value = 2399.111
$ws.Range("K77").Value = 7334.7145
$ws.Range("L77").Value = 11995.555
$ws.Range("M77").Value = -2966.7145
$ws.Range("N77").Value = -20731.555
$ws.Range("H132").Value = 1804.0869
$ws.Range("I132").Value = 980.5
$ws.Range("K132").Value = 2941.5
$ws.Range("M132").Value = -411.5
$ws.Range("H136").Value = 1207.2703
$ws.Range("I136").Value = 919.8148
$ws.Range("J136").Value = 1983.4
$ws.Range("K136").Value = 2759.4444
$ws.Range("L136").Value = 5950.200000000001
$ws.Range("M136").Value = -209.4443999999999
$ws.Range("N136").Value = -11050.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 904.3333
$ws.Range("I22").Value = 904.3333
$ws.Range("K22").Value = 904.3333
$ws.Range("M22").Value = -731.3333
$ws.Range("H40").Value = 18000
$ws.Range("J40").Value = 18000
$ws.Range("L40").Value = 18000
$ws.Range("N40").Value = -18530
$ws.Range("H96").Value = 20000
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("H134").Value = 2848.2778
$ws.Range("I134").Value = 1841.4445
$ws.Range("J134").Value = 3452.3777
$ws.Range("K134").Value = 5524.333500000001
$ws.Range("L134").Value = 10357.1331
$ws.Range("M134").Value = -2989.333500000001
$ws.Range("N134").Value = -15427.1331

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1505.2222
$ws.Range("I22").Value = 374.5
$ws.Range("J22").Value = 3766.6667
$ws.Range("K22").Value = 374.5
$ws.Range("L22").Value = 3766.6667
$ws.Range("M22").Value = -24.5
$ws.Range("N22").Value = -4466.6667
$ws.Range("H31").Value = 4419.11
$ws.Range("I31").Value = 2882.7856
$ws.Range("J31").Value = 4698.4414
$ws.Range("K31").Value = 2882.7856
$ws.Range("L31").Value = 4698.4414
$ws.Range("M31").Value = -2587.7856
$ws.Range("N31").Value = -5288.4414
$ws.Range("H34").Value = 4419.11
$ws.Range("I34").Value = 2882.7856
$ws.Range("J34").Value = 4698.4414
$ws.Range("K34").Value = 2882.7856
$ws.Range("L34").Value = 4698.4414
$ws.Range("M34").Value = -2680.7856
$ws.Range("N34").Value = -5102.4414
$ws.Range("H58").Value = 1657.0968
$ws.Range("I58").Value = 1299.1072
$ws.Range("J58").Value = 4998.3335
$ws.Range("K58").Value = 1299.1072
$ws.Range("L58").Value = 4998.3335
$ws.Range("M58").Value = -1096.1072
$ws.Range("N58").Value = -5404.3335
$ws.Range("H132").Value = 56657.69
$ws.Range("I132").Value = 1370.0588
$ws.Range("J132").Value = 161089.89
$ws.Range("K132").Value = 4110.1764
$ws.Range("L132").Value = 483269.67
$ws.Range("M132").Value = -1580.1764
$ws.Range("N132").Value = -488329.67
$ws.Range("H134").Value = 342850.56
$ws.Range("I134").Value = 1037.1177
$ws.Range("J134").Value = 2003087.2
$ws.Range("K134").Value = 3111.3531
$ws.Range("L134").Value = 6009261.6
$ws.Range("M134").Value = -576.3531000000003
$ws.Range("N134").Value = -6014331.6
$ws.Range("H136").Value = 1657.0968
$ws.Range("I136").Value = 1299.1072
$ws.Range("J136").Value = 4998.3335
$ws.Range("K136").Value = 3897.3216
$ws.Range("L136").Value = 14995.0005
$ws.Range("M136").Value = -1347.3216
$ws.Range("N136").Value = -20095.0005

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 24268358
$ws.Range("I33").Value = 200
$ws.Range("J33").Value = 29661282
$ws.Range("K33").Value = 1200
$ws.Range("L33").Value = 177967692
$ws.Range("M33").Value = -917
$ws.Range("N33").Value = -177968258
$ws.Range("H60").Value = 523.3333
$ws.Range("I60").Value = 328
$ws.Range("J60").Value = 1500
$ws.Range("K60").Value = 984
$ws.Range("L60").Value = 4500
$ws.Range("M60").Value = -733
$ws.Range("N60").Value = -5002
$ws.Range("H75").Value = 2690.75
$ws.Range("I75").Value = 350.6
$ws.Range("J75").Value = 3754.4546
$ws.Range("K75").Value = 1051.8
$ws.Range("L75").Value = 11263.3638
$ws.Range("M75").Value = -53.80000000000018
$ws.Range("N75").Value = -13259.3638
$ws.Range("H78").Value = 2690.75
$ws.Range("I78").Value = 350.6
$ws.Range("J78").Value = 3754.4546
$ws.Range("K78").Value = 3155.4
$ws.Range("L78").Value = 33790.0914
$ws.Range("M78").Value = 1836.6
$ws.Range("N78").Value = -43774.0914
$ws.Range("H109").Value = 1900.8846
$ws.Range("I109").Value = 920.86664
$ws.Range("K109").Value = 2762.59992
$ws.Range("M109").Value = -1722.59992
$ws.Range("H113").Value = 2845.1738
$ws.Range("I113").Value = 3906.4
$ws.Range("J113").Value = 855.375
$ws.Range("K113").Value = 11719.2
$ws.Range("L113").Value = 2566.125
$ws.Range("M113").Value = -9549.200000000001
$ws.Range("N113").Value = -6906.125
$ws.Range("H131").Value = 3438.1226
$ws.Range("J131").Value = 1581.4634
$ws.Range("L131").Value = 4744.3902
$ws.Range("N131").Value = -14824.3902
$ws.Range("H137").Value = 27034118
$ws.Range("I137").Value = 3070
$ws.Range("J137").Value = 41675936
$ws.Range("K137").Value = 9210
$ws.Range("L137").Value = 125027808
$ws.Range("M137").Value = -4110
$ws.Range("N137").Value = -125038008

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3636.2222
$ws.Range("I80").Value = 3922.923
$ws.Range("J80").Value = 3370
$ws.Range("K80").Value = 3922.923
$ws.Range("L80").Value = 3370
$ws.Range("M80").Value = -2924.923
$ws.Range("N80").Value = -5366
$ws.Range("H83").Value = 3636.2222
$ws.Range("I83").Value = 3922.923
$ws.Range("J83").Value = 3370
$ws.Range("K83").Value = 19614.615
$ws.Range("L83").Value = 16850
$ws.Range("M83").Value = -14622.615
$ws.Range("N83").Value = -26834
$ws.Range("H122").Value = 1123.6
$ws.Range("I122").Value = 1136.4445
$ws.Range("J122").Value = 1008
$ws.Range("K122").Value = 3409.3335
$ws.Range("L122").Value = 3024
$ws.Range("M122").Value = -959.3335000000002
$ws.Range("N122").Value = -7924
$ws.Range("H132").Value = 2406.8918
$ws.Range("I132").Value = 1377.0952
$ws.Range("K132").Value = 4131.2856
$ws.Range("M132").Value = -1601.2856

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6072.857
$ws.Range("I46").Value = 1174.2858
$ws.Range("J46").Value = 10971.429
$ws.Range("K46").Value = 1174.2858
$ws.Range("L46").Value = 10971.429
$ws.Range("M46").Value = -986.2858000000001
$ws.Range("N46").Value = -11347.429
$ws.Range("H132").Value = 1905.7205
$ws.Range("I132").Value = 1397.6865
$ws.Range("J132").Value = 3214.8845
$ws.Range("K132").Value = 4193.0595
$ws.Range("L132").Value = 9644.6535
$ws.Range("M132").Value = -1663.0595
$ws.Range("N132").Value = -14704.6535
$ws.Range("H136").Value = 1963.2759
$ws.Range("I136").Value = 1611.85
$ws.Range("K136").Value = 4835.549999999999
$ws.Range("M136").Value = -2285.549999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1553.1777
$ws.Range("I132").Value = 1185.0588
$ws.Range("J132").Value = 2691
$ws.Range("K132").Value = 3555.1764
$ws.Range("L132").Value = 8073
$ws.Range("M132").Value = -1025.1764
$ws.Range("N132").Value = -13133
$ws.Range("H136").Value = 286519.16
$ws.Range("I136").Value = 370940.22
$ws.Range("J136").Value = 1598.125
$ws.Range("K136").Value = 1112820.66
$ws.Range("L136").Value = 4794.375
$ws.Range("M136").Value = -1110270.66
$ws.Range("N136").Value = -9894.375

# --- Clear M96 on BSM (cell removed in target) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M96").ClearContents()
